$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.83279933333333
$ws.Range("H2").Value = 50.498398
$ws.Range("I2").Value = 0.04383102208811961
$ws.Range("J2").Value = 0.04383102208811961
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 148.4761865521018
$ws.Range("R2").Value = 1336.285678968916
$ws.Range("S2").Value = 0.002811824349710392
$ws.Range("T2").Value = 0.002811824349710392
$ws.Range("G3").Value = 16.83279933333333
$ws.Range("H3").Value = 50.498398
$ws.Range("I3").Value = 0.04383102208811961
$ws.Range("J3").Value = 0.04383102208811961
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 920.9482898736434
$ws.Range("R3").Value = 8288.53460886279
$ws.Range("S3").Value = 0.01744080910498167
$ws.Range("T3").Value = 0.01744080910498168
$ws.Range("G4").Value = 16.83279933333333
$ws.Range("H4").Value = 50.498398
$ws.Range("I4").Value = 0.04383102208811961
$ws.Range("J4").Value = 0.04383102208811961
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 368.7757171518911
$ws.Range("R4").Value = 3318.98145436702
$ws.Range("S4").Value = 0.006983830640785819
$ws.Range("T4").Value = 0.006983830640785819
$ws.Range("G5").Value = 16.83279933333333
$ws.Range("H5").Value = 50.498398
$ws.Range("I5").Value = 0.04383102208811961
$ws.Range("J5").Value = 0.04383102208811961
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 876.2626614706273
$ws.Range("R5").Value = 7886.363953235646
$ws.Range("S5").Value = 0.01659455799264172
$ws.Range("T5").Value = 0.01659455799264172
$ws.Range("I6").Value = 0.8654671168650652
$ws.Range("J6").Value = 0.8654671168650654
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 2931.742199395285
$ws.Range("R6").Value = 26385.67979455757
$ws.Range("S6").Value = 0.05552098484453208
$ws.Range("T6").Value = 0.05552098484453209
$ws.Range("I7").Value = 0.8654671168650652
$ws.Range("J7").Value = 0.8654671168650654
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.3443781607815578
$ws.Range("T7").Value = 0.344378160781558
$ws.Range("I8").Value = 0.8654671168650652
$ws.Range("J8").Value = 0.8654671168650654
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 7281.674975583176
$ws.Range("R8").Value = 65535.07478024859
$ws.Range("S8").Value = 0.1378994940433548
$ws.Range("T8").Value = 0.1378994940433548
$ws.Range("I9").Value = 0.8654671168650652
$ws.Range("J9").Value = 0.8654671168650654
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 17302.27777291669
$ws.Range("R9").Value = 155720.4999562502
$ws.Range("S9").Value = 0.3276684771956204
$ws.Range("T9").Value = 0.3276684771956205
$ws.Range("G10").Value = 34.50825133333333
$ws.Range("H10").Value = 103.524754
$ws.Range("I10").Value = 0.08985623225594501
$ws.Range("J10").Value = 0.08985623225594502
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 304.3851151013631
$ws.Range("R10").Value = 2739.466035912268
$ws.Range("S10").Value = 0.005764409082739185
$ws.Range("T10").Value = 0.005764409082739186
$ws.Range("G11").Value = 34.50825133333333
$ws.Range("H11").Value = 103.524754
$ws.Range("I11").Value = 0.08985623225594501
$ws.Range("J11").Value = 0.08985623225594502
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 1887.999400612463
$ws.Range("R11").Value = 16991.99460551217
$ws.Range("S11").Value = 0.03575470794448149
$ws.Range("T11").Value = 0.0357547079444815
$ws.Range("G12").Value = 34.50825133333333
$ws.Range("H12").Value = 103.524754
$ws.Range("I12").Value = 0.08985623225594501
$ws.Range("J12").Value = 0.08985623225594502
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 756.0124065583844
$ws.Range("R12").Value = 6804.111659025461
$ws.Range("S12").Value = 0.0143172729769569
$ws.Range("T12").Value = 0.0143172729769569
$ws.Range("G13").Value = 34.50825133333333
$ws.Range("H13").Value = 103.524754
$ws.Range("I13").Value = 0.08985623225594501
$ws.Range("J13").Value = 0.08985623225594502
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 1796.391173995895
$ws.Range("R13").Value = 16167.52056596306
$ws.Range("S13").Value = 0.03401984225176743
$ws.Range("T13").Value = 0.03401984225176744
$ws.Range("G14").Value = 0.324754
$ws.Range("H14").Value = 0.974262
$ws.Range("I14").Value = 0.0008456287908700705
$ws.Range("J14").Value = 0.0008456287908700706
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 2.864540504089333
$ws.Range("R14").Value = 25.780864536804
$ws.Range("S14").Value = [double]"5.424832713698256E-05"
$ws.Range("T14").Value = [double]"5.424832713698256E-05"
$ws.Range("G15").Value = 0.324754
$ws.Range("H15").Value = 0.974262
$ws.Range("I15").Value = 0.0008456287908700705
$ws.Range("J15").Value = 0.0008456287908700706
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 17.76778983739
$ws.Range("R15").Value = 159.91010853651
$ws.Range("S15").Value = 0.0003364842892686943
$ws.Range("T15").Value = 0.0003364842892686944
$ws.Range("G16").Value = 0.324754
$ws.Range("H16").Value = 0.974262
$ws.Range("I16").Value = 0.0008456287908700705
$ws.Range("J16").Value = 0.0008456287908700706
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 7.114763675153333
$ws.Range("R16").Value = 64.03287307638
$ws.Range("S16").Value = 0.0001347385477011226
$ws.Range("T16").Value = 0.0001347385477011226
$ws.Range("G17").Value = 0.324754
$ws.Range("H17").Value = 0.974262
$ws.Range("I17").Value = 0.0008456287908700705
$ws.Range("J17").Value = 0.0008456287908700706
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 16.905673187686
$ws.Range("R17").Value = 152.151058689174
$ws.Range("S17").Value = 0.0003201576267632709
$ws.Range("T17").Value = 0.000320157626763271

Write-Host "Applied 174 cell updates"
